$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared / rich text edits -------------------------------------------------
# A8 holds the rich "Volume 31   Number  5" string; only the trailing issue
# number run ("5" -> "6") changes.
$volRange = $ws.Range("A8")
$volText = $volRange.Characters().Text
$volRange.Characters($volText.Length, 1).Text = "6"

# C9 holds the rich "Report Covering the Week  1/29/2024  Through  2/4/2024"
# string; the two date runs change. Edit the right-most one first so the
# left edit's length change doesn't shift the already-computed offset.
$dateRange = $ws.Range("C9")
$dateRange.Characters(47, 8).Text = "2/11/2024"
$dateRange.Characters(27, 9).Text = "2/5/2024"

# --- Numeric grid updates (rows 14-30, columns C-N) ---------------------------
$updates = @(
  @(14, 3, 3),
  @(14, 4, 4),
  @(14, 5, -25),
  @(14, 6, 8),
  @(14, 7, 12),
  @(14, 8, -33.333333333333),
  @(14, 9, 11),
  @(14, 10, 15),
  @(14, 11, -26.666666666666),
  @(14, 12, 0),
  @(14, 13, 37.5),
  @(14, 14, -79.629629629629),
  @(15, 3, 9),
  @(15, 4, 11),
  @(15, 5, -18.181818181818),
  @(15, 6, 35),
  @(15, 7, 28),
  @(15, 8, 25),
  @(15, 9, 46),
  @(15, 10, 47),
  @(15, 11, -2.127659574468),
  @(15, 12, -14.814814814814),
  @(15, 13, 39.393939393939),
  @(15, 14, -24.590163934426),
  @(16, 3, 74),
  @(16, 4, 74),
  @(16, 5, 0),
  @(16, 6, 359),
  @(16, 7, 314),
  @(16, 8, 14.331210191082),
  @(16, 9, 582),
  @(16, 10, 508),
  @(16, 11, 14.566929133858),
  @(16, 12, 22.268907563025),
  @(16, 13, 22.012578616352),
  @(16, 14, -70.21494370522),
  @(17, 3, 130),
  @(17, 4, 137),
  @(17, 5, -5.109489051094),
  @(17, 6, 529),
  @(17, 7, 519),
  @(17, 8, 1.926782273603),
  @(17, 9, 807),
  @(17, 10, 781),
  @(17, 11, 3.329065300896),
  @(17, 12, 16.115107913669),
  @(17, 13, 83.40909090909),
  @(17, 14, -3.699284009546),
  @(18, 3, 54),
  @(18, 4, 66),
  @(18, 5, -18.181818181818),
  @(18, 6, 211),
  @(18, 7, 231),
  @(18, 8, -8.658008658008),
  @(18, 9, 306),
  @(18, 10, 352),
  @(18, 11, -13.068181818181),
  @(18, 12, -4.075235109717),
  @(18, 13, -21.739130434782),
  @(18, 14, -85.982592762253),
  @(19, 3, 180),
  @(19, 4, 141),
  @(19, 5, 27.659574468085),
  @(19, 6, 733),
  @(19, 7, 537),
  @(19, 8, 36.499068901303),
  @(19, 9, 1064),
  @(19, 10, 815),
  @(19, 11, 30.552147239263),
  @(19, 12, 24.29906542056),
  @(19, 13, 139.63963963964),
  @(19, 14, 35.714285714285),
  @(20, 3, 76),
  @(20, 4, 121),
  @(20, 5, -37.190082644628),
  @(20, 6, 313),
  @(20, 7, 411),
  @(20, 8, -23.844282238442),
  @(20, 9, 472),
  @(20, 10, 636),
  @(20, 11, -25.786163522012),
  @(20, 12, -13.553113553113),
  @(20, 13, 115.525114155251),
  @(20, 14, -73.965802537231),
  @(21, 3, 526),
  @(21, 4, 554),
  @(21, 5, -5.054151624548),
  @(21, 6, 2188),
  @(21, 7, 2052),
  @(21, 8, 6.62768031189),
  @(21, 9, 3288),
  @(21, 10, 3154),
  @(21, 11, 4.248573240329),
  @(21, 12, 11.193777477172),
  @(21, 13, 63.419483101391),
  @(21, 14, -57.22648627553),
  @(22, 3, 4),
  @(22, 4, 7),
  @(22, 5, -42.857142857142),
  @(22, 6, 27),
  @(22, 7, 15),
  @(22, 8, 80),
  @(22, 9, 42),
  @(22, 10, 22),
  @(22, 11, 90.90909090909),
  @(22, 12, 2.439024390243),
  @(22, 13, 35.483870967741),
  @(23, 3, 24),
  @(23, 4, 35),
  @(23, 5, -31.428571428571),
  @(23, 6, 120),
  @(23, 7, 128),
  @(23, 8, -6.25),
  @(23, 9, 186),
  @(23, 10, 202),
  @(23, 11, -7.920792079207),
  @(23, 12, 12.727272727272),
  @(23, 13, 75.471698113207),
  @(24, 3, 288),
  @(24, 4, 297),
  @(24, 5, -3.030303030303),
  @(24, 6, 1289),
  @(24, 7, 1247),
  @(24, 8, 3.36808340016),
  @(24, 9, 1892),
  @(24, 10, 1834),
  @(24, 11, 3.162486368593),
  @(24, 12, 7.622298065984),
  @(24, 13, 39.015429831006),
  @(25, 3, 191),
  @(25, 4, 193),
  @(25, 5, -1.036269430051),
  @(25, 6, 776),
  @(25, 7, 743),
  @(25, 8, 4.441453566621),
  @(25, 9, 1156),
  @(25, 10, 1066),
  @(25, 11, 8.442776735459),
  @(25, 12, 14.682539682539),
  @(25, 13, 10.30534351145),
  @(26, 3, 10),
  @(26, 4, 15),
  @(26, 5, -33.333333333333),
  @(26, 6, 54),
  @(26, 7, 48),
  @(26, 8, 12.5),
  @(26, 9, 70),
  @(26, 10, 76),
  @(26, 11, -7.894736842105),
  @(26, 12, -18.60465116279),
  @(27, 3, 19),
  @(27, 4, 18),
  @(27, 5, 5.555555555555),
  @(27, 6, 76),
  @(27, 7, 75),
  @(27, 8, 1.333333333333),
  @(27, 9, 114),
  @(27, 10, 125),
  @(27, 11, -8.8),
  @(27, 12, 31.03448275862),
  @(28, 3, 11),
  @(28, 4, 10),
  @(28, 5, 10),
  @(28, 6, 27),
  @(28, 7, 31),
  @(28, 8, -12.903225806451),
  @(28, 9, 40),
  @(28, 10, 37),
  @(28, 11, 8.108108108108),
  @(28, 12, -27.272727272727),
  @(28, 13, 0),
  @(28, 14, -66.94214876033),
  @(29, 3, 7),
  @(29, 4, 6),
  @(29, 5, 16.666666666666),
  @(29, 6, 20),
  @(29, 7, 22),
  @(29, 8, -9.090909090909),
  @(29, 9, 31),
  @(29, 10, 28),
  @(29, 11, 10.714285714285),
  @(29, 12, -39.215686274509),
  @(29, 13, -16.216216216216),
  @(29, 14, -72.321428571428),
  @(30, 3, 1),
  @(30, 6, 1),
  @(30, 7, 2),
  @(30, 8, -50),
  @(30, 9, 2),
  @(30, 11, -50),
  @(30, 12, -66.666666666666)
)
foreach ($u in $updates) {
  $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# --- C30 / F30 change from "n/a" text to an actual number; give them the
#     same numeric style (#,##0) the rest of the C/F column already uses.
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("F30").NumberFormat = "#,##0"
